$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 38
# from 45204 (2023-10-05) to 45207 (2023-10-08)
$ws.Range("C2:C38").Value = 45207
